# Add a new list item ("Handle state changes and populate states") right
# after the "Apply logic for special cases" bullet, keeping the same
# ListParagraph style / numbering (numId 9, ilvl 0) and splitting the new
# text across two runs, as in the authored change.

$d = $word.ActiveDocument

# Locate the "Apply logic for special cases" paragraph (currently the
# last paragraph in the document body).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.TrimEnd("`r") -eq "Apply logic for special cases") {
        $targetPara = $candidate
    }
}

# Insert a new paragraph right after it; it inherits the paragraph
# formatting (style + numbering) of $targetPara, matching the diff.
[void]$targetPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)

# Insert the new content as two distinct runs via raw OOXML so the run
# boundary ("Handle state changes" / " and populate states") is preserved
# exactly like the authored edit, instead of Word's usual auto-merge of
# adjacent same-formatted runs.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr>' +
        '<w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr>' +
    '</w:pPr>' +
    '<w:r><w:t>Handle state changes</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and populate states</w:t></w:r>' +
'</w:p>'

[void]$newPara.Range.InsertXML($newParaXml)
